# reviewdb.xlsx — "Add files via upload"
#
# The sheet is a flat review table (appid, keyword, email, recovery, time,
# review, blue) that runs from row 2 through row 18. The upload added one
# more review row (appended as row 19, formatted like the other data rows)
# and dropped a stray leftover value that had been sitting by itself out at
# I27 (far outside the table, column I / row 27 — not part of the A:G data
# block at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the orphan cell at I27 (the whole row only has that one cell, so
# removing the row is equivalent to removing the cell; it also pulls the
# used range back in from column I to column G).
$ws.Rows("27").Delete()

# New row 19 mirrors the look of the other "bitcoin / com.hamxa.shaynachim"
# rows (e.g. row 5): bold appid font in A, centered/wrapped email-style
# formatting in C:D. Copy that formatting down first, then overwrite with
# the new row's actual values.
$ws.Range("A5:G5").Copy($ws.Range("A19:G19"))

$ws.Range("A19").Value = "com.hamxa.shaynachim"
$ws.Range("B19").Value = "bitcoin"
$ws.Range("C19").Value = "ctamar115@gmail.com"
$ws.Range("D19").Value = "nirh94846@gmail.com"
$ws.Range("E19").Value = "27/5/2019 15:59"
$ws.Range("F19").Value = "the game is just realistic and resembles the lifetime we are in today.it is also adventurous.teaches us how to save money and make budjets for future use"
$ws.Range("G19").Value = "no"

$ws.Rows("19").RowHeight = 13.8

# Match the new selection / scroll position left in the worksheet view.
$ws.Activate()
$ws.Range("C19:D19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 3
